# Apply the column-insertion edit to ODP Site 516 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank columns at their respective positions.
# Process right-to-left (using the ORIGINAL/old column letters) so that each
# insertion does not disturb the addressing of columns still to be processed.
$ws.Columns("S").Insert()   # new column before old S ("MgCa Temperature anomaly_Original - Coretop")
$ws.Columns("R").Insert()   # new column before old R ("MgCa Temperature anomaly_BAYMAG - ERSST")
$ws.Columns("M").Insert()   # new column before old M ("MgCa Temperature_BAYMAG")

# --- New column M: MgCa Temperature_Original ---
$ws.Range("M1").Value = "MgCa Temperature_Original"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108  # xlCenter, matches other header cells
$ws.Range("M2").Value = 24.789288925

# --- New column S: MgCa Temperature anomaly_Original - ERSST ---
$ws.Range("S1").Value = "MgCa Temperature anomaly_Original - ERSST"
$ws.Range("S1").Font.Bold = $true
$ws.Range("S1").HorizontalAlignment = -4108
$ws.Range("S2").Value = 4.049772560796445

# --- New column U: MgCa Coretop modelled temperature ---
$ws.Range("U1").Value = "MgCa Coretop modelled temperature"
$ws.Range("U1").Font.Bold = $true
$ws.Range("U1").HorizontalAlignment = -4108
$ws.Range("U2").Value = 22.0667

# --- Updated values in shifted columns ---
# R2 (was Q2, ERSST_V5 Temperature) changed from 20.43029746 to 20.74
$ws.Range("R2").Value = 20.74

# T2 (was R2, MgCa Temperature anomaly_BAYMAG - ERSST) changed value
$ws.Range("T2").Value = 5.313683635796444

# W2 (was T2, MgCa Temperature anomaly_BAYMAG - Coretop) changed precision
$ws.Range("W2").Value = 3.986499999999999

Write-Host ("Final dimension: " + $ws.UsedRange.Address())
